$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1487.511676915803
$ws.Range("D2").Value = 1883.468544629751
$ws.Range("B3").Value = 2597.1443555309761
$ws.Range("D3").Value = 2631.0266361136132
$ws.Range("B4").Value = 3749.2276811183551
$ws.Range("D4").Value = 3361.049035746008
$ws.Range("B5").Value = 4359.7024365325342
$ws.Range("D5").Value = 3702.6385577445481
$ws.Range("B6").Value = 4867.7790804145307
$ws.Range("D6").Value = 3936.7122999905491
$ws.Range("B7").Value = 5477.0020110422483
$ws.Range("D7").Value = 4019.2463258656821
$ws.Range("B8").Value = 5870.9757057273737
$ws.Range("D8").Value = 4056.8255077806839
$ws.Range("B9").Value = 6108.2672100881236
$ws.Range("D9").Value = 4019.119389661651
$ws.Range("B10").Value = 6509.2126481888808
$ws.Range("D10").Value = 3938.0230635159801
$ws.Range("B11").Value = 6656.6219839515343
$ws.Range("D11").Value = 3898.0037236374642
$ws.Range("B12").Value = 6830.4194273670109
$ws.Range("D12").Value = 3764.5337586066212
$ws.Range("B13").Value = 7114.3915848310071
$ws.Range("D13").Value = 3530.0264109258378
$ws.Range("B14").Value = 7253.2620127701339
$ws.Range("D14").Value = 3308.580783525048
$ws.Range("B15").Value = 7389.2408480578224
$ws.Range("D15").Value = 3141.431113705878
$ws.Range("B16").Value = 7599.0573863269792
$ws.Range("D16").Value = 2867.4796875735278
$ws.Range("B17").Value = 7726.4113143407349
$ws.Range("D17").Value = 2554.1921565376101
$ws.Range("B18").Value = 7847.7566768792394
$ws.Range("D18").Value = 2191.9683281496891
$ws.Range("B19").Value = 7934.4090808749643
$ws.Range("D19").Value = 1888.107444918026
$ws.Range("B20").Value = 8091.6896000928291
$ws.Range("D20").Value = 1502.7110910298129
$ws.Range("B21").Value = 8128.6249263326026
$ws.Range("D21").Value = 1351.854507810848
$ws.Range("B22").Value = 8249.5272012321475
$ws.Range("D22").Value = 865.92133840142242
$ws.Range("B23").Value = 8271.349737002558
$ws.Range("D23").Value = 594.50865899265602
$ws.Range("B24").Value = 8283.3270439147764
$ws.Range("D24").Value = 444.95209229970538
$ws.Range("B25").Value = 8324.5893621095711
$ws.Range("D25").Value = 103.010632640355
$ws.Range("B26").Value = 8321.1841228069097
$ws.Range("D26").Value = 103.0106324305509
$ws.Range("B27").Value = 8321.2929399939439
$ws.Range("D27").Value = 103.0106324113675
$ws.Range("B28").Value = 8321.2929399939439
$ws.Range("D28").Value = 103.0106324113675
$ws.Range("B29").Value = 8321.2929399939439
$ws.Range("D29").Value = 103.0106324113675
$ws.Range("B30").Value = 8321.2929399939439
$ws.Range("D30").Value = 103.0106324113675
$ws.Range("B31").Value = 8321.2929399939439
$ws.Range("D31").Value = 103.0106324113675
$ws.Range("B32").Value = 8321.2929399939439
$ws.Range("D32").Value = 103.0106324113675
$ws.Range("B33").Value = 8321.2929399939439
$ws.Range("D33").Value = 103.0106324113675
$ws.Range("B34").Value = 8321.2929399939439
$ws.Range("D34").Value = 103.0106324113675
$ws.Range("B35").Value = 8321.2929399939439
$ws.Range("D35").Value = 103.0106324113675
$ws.Range("B36").Value = 8321.2929399939439
$ws.Range("D36").Value = 103.0106324113675
$ws.Range("B37").Value = 8321.2929399939439
$ws.Range("D37").Value = 103.0106324113675
$ws.Range("B38").Value = 8321.2929399939439
$ws.Range("D38").Value = 103.0106324113675
$ws.Range("B39").Value = 8321.2929399939439
$ws.Range("D39").Value = 103.0106324113675
$ws.Range("B40").Value = 8321.2929399939439
$ws.Range("D40").Value = 103.0106324113675
$ws.Range("B41").Value = 8321.2929399939439
$ws.Range("D41").Value = 103.0106324113675
$ws.Range("B42").Value = 8321.2929399939439
$ws.Range("D42").Value = 103.0106324113675
$ws.Range("B43").Value = 8321.2929399939439
$ws.Range("D43").Value = 103.0106324113675
$ws.Range("B44").Value = 8321.2929399939439
$ws.Range("D44").Value = 103.0106324113675
$ws.Range("B45").Value = 8321.2929399939439
$ws.Range("D45").Value = 103.0106324113675
$ws.Range("B46").Value = 8321.2929399939439
$ws.Range("D46").Value = 103.0106324113675
$ws.Range("B47").Value = 8321.2929399939439
$ws.Range("D47").Value = 103.0106324113675
$ws.Range("B48").Value = 8321.2929399939439
$ws.Range("D48").Value = 103.0106324113675
$ws.Range("B49").Value = 8321.2929399939439
$ws.Range("D49").Value = 103.0106324113675
$ws.Range("B50").Value = 8321.2929399939439
$ws.Range("D50").Value = 103.0106324113675
$ws.Range("B51").Value = 8321.2929399939439
$ws.Range("D51").Value = 103.0106324113675
$ws.Range("B52").Value = 8321.2929399939439
$ws.Range("D52").Value = 103.0106324113675
$ws.Range("B53").Value = 8321.2929399939439
$ws.Range("D53").Value = 103.0106324113675
$ws.Range("B54").Value = 8321.2929399939439
$ws.Range("D54").Value = 103.0106324113675
$ws.Range("B55").Value = 8321.2929399939439
$ws.Range("D55").Value = 103.0106324113675
$ws.Range("B56").Value = 8321.2929399939439
$ws.Range("D56").Value = 103.0106324113675
$ws.Range("B57").Value = 8321.2929399939439
$ws.Range("D57").Value = 103.0106324113675
$ws.Range("B58").Value = 8321.2929399939439
$ws.Range("D58").Value = 103.0106324113675
$ws.Range("B59").Value = 8321.2929399939439
$ws.Range("D59").Value = 103.0106324113675
$ws.Range("B60").Value = 8321.2929399939439
$ws.Range("D60").Value = 103.0106324113675
$ws.Range("B61").Value = 8321.2929399939439
$ws.Range("D61").Value = 103.0106324113675
$ws.Range("B62").Value = 8321.2929399939439
$ws.Range("D62").Value = 103.0106324113675

$ws.Range("F8").Select()